# Update Work Week and Social Spending
# (Actually: refresh the Slovakia GDP per Capita "Data" sheet with a new
# 1985-2016 series, replacing the old 1990-2010 series.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

function Set-TextValue($range, [string]$val) {
    # Force the cell to be written as a shared string (t="s"), even though
    # the text looks like a number, then drop the temporary number-format
    # override so no extra cell style sticks around.
    $range.NumberFormat = "@"
    $range.Value = $val
    $range.ClearFormats()
}

$rows = @(
    @{Year=1985; Data="11815"},
    @{Year=1986; Data="12242"},
    @{Year=1987; Data="12473"},
    @{Year=1988; Data="12659"},
    @{Year=1989; Data="12734"},
    @{Year=1990; Data="12374"},
    @{Year=1991; Data="10537.4103281024"},
    @{Year=1992; Data="9838.61485840115"},
    @{Year=1993; Data="10481.1152836241"},
    @{Year=1994; Data="11063.9487142439"},
    @{Year=1995; Data="11874.3980037823"},
    @{Year=1996; Data="12622.9095750226"},
    @{Year=1997; Data="13335.4589939449"},
    @{Year=1998; Data="13822.1874047156"},
    @{Year=1999; Data="13750.3330903464"},
    @{Year=2000; Data="13904.9853726498"},
    @{Year=2001; Data="14361.9012127386"},
    @{Year=2002; Data="14984.7311249787"},
    @{Year=2003; Data="15773.0947996538"},
    @{Year=2004; Data="16570.8618948388"},
    @{Year=2005; Data="17649.5203820258"},
    @{Year=2006; Data="19099.4276773689"},
    @{Year=2007; Data="21109.9388472468"},
    @{Year=2008; Data="22231.9881558054"},
    @{Year=2009; Data="20953.0366618401"},
    @{Year=2010; Data="21941.2121615904"},
    @{Year=2011; Data="22483"},
    @{Year=2012; Data="22816"},
    @{Year=2013; Data="23132"},
    @{Year=2014; Data="23703"},
    @{Year=2015; Data="24588"},
    @{Year=2016; Data="25364"}
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $item = $rows[$i]

    $ws.Range("A$r").Value = 703
    $ws.Range("B$r").Value = "Slovakia"
    $ws.Range("C$r").Value = "GDP per Capita"
    $ws.Range("D$r").Value = $item.Year
    Set-TextValue $ws.Range("E$r") $item.Data
}
